# Updated cryptos list values (Price and Volume(1h) columns) per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.523.06"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.36%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.842.45"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.12%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -1.35%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "333.50"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4631"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3848"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.57%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.02"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07907"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9927"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.99%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.85%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.850.51"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.926"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.111"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.008"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.75"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06679"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.39%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.04%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.005"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.538.22"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.379"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.59%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.311"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.92"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.56%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.104"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.397"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "119.82"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.9742"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.98%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09382"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.79%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.274"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.90%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.339"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.45%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06025"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.42%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02227"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.310"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.19%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.99%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5882"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.22%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.43%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "10.29"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.240"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.45%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5574"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.13"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.86%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.06695"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.46%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "110.97"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.005"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "70.07"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.88%  "
